$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("YDS")

$ws.Range("B2").Value = "NCT(2.279906628067043, 1.519524995518351, -0.277711842643387, 2.179763599604267)"
$ws.Range("C2").Value = "JSU(-1.3779226476176363, 1.2330460899897009, 2.104395625237625, 4.884614310058583)"
$ws.Range("D2").Value = "NIG(0.8437896977083359, 0.5779014779080955, 1.5786788790712465, 3.381658367130838)"
$ws.Range("E2").Value = "NIG(1.0559770899709198, 0.7830990314371007, 4.72437133049091, 6.774330604560702)"
